# "Generate Report for Handoff"
#
# The CI job that produces this localization-status report re-ran the
# handoff step for the "474c97d6-ac65-40bf-b1bd-69f6d19a0c9a" source file,
# recording a newer "Latest Handoff Datetime" for both locale sheets
# (row 6 in each of the zh-cn and de-de tables, column D).
#
# zh-cn: 2016-03-09 04:48:25 -> 2016-03-09 04:48:59
# de-de: 2016-03-09 04:48:27 -> 2016-03-09 04:49:02

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 04:48:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 04:49:02"
